$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.308.88'
$ws.Range("E2").Value = '  +0.54%  '
$ws.Range("D3").Value = '2.374.81'
$ws.Range("E3").Value = '  +0.49%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '552.67'
$ws.Range("E5").Value = '  +2.06%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.45'
$ws.Range("E6").Value = '  +1.08%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -0.23%  '
$ws.Range("D9").Value = '2.376.46'
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.107'
$ws.Range("E10").Value = '  +2.46%  '
$ws.Range("E11").Value = '  +1.96%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.36'
$ws.Range("E12").Value = '  +2.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.350'
$ws.Range("E13").Value = '  +2.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.45'
$ws.Range("E14").Value = '  +2.19%  '
$ws.Range("E15").Value = '  +4.84%  '
$ws.Range("D16").Value = '2.803.16'
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("D17").Value = '61.246.63'
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").Value = '2.372.57'
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.99'
$ws.Range("E19").Value = '  +3.40%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.15'
$ws.Range("E20").Value = '  +1.69%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '320.09'
$ws.Range("E21").Value = '  +1.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.68'
$ws.Range("E22").Value = '  +1.32%  '
$ws.Range("E23").Value = '  +0.06%  '
$ws.Range("E24").Value = '  -8.62%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '64.19'
$ws.Range("E25").Value = '  +1.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.86'
$ws.Range("E26").Value = '  +4.91%  '
$ws.Range("D28").Value = '2.492.03'
$ws.Range("E28").Value = '  +0.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.15'
$ws.Range("E29").Value = '  +2.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '516.65'
$ws.Range("E30").Value = '  +1.62%  '
$ws.Range("D31").Value = '0.0₃0899'
$ws.Range("E31").Value = '  +0.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.39'
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("E33").Value = '  +2.71%  '
$ws.Range("E34").Value = '  +2.92%  '
$ws.Range("E35").Value = '  -1.93%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.51'
$ws.Range("E37").Value = '  +4.63%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.68'
$ws.Range("E38").Value = '  +1.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.88'
$ws.Range("E39").Value = '  +5.14%  '
$ws.Range("E40").Value = '  +1.51%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '18.47'
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '147.06'
$ws.Range("E42").Value = '  +6.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.999'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.34'
$ws.Range("E44").Value = '  +3.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '147.28'
$ws.Range("E45").Value = '  +6.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.14'
$ws.Range("E46").Value = '  -0.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.59'
$ws.Range("E47").Value = '  +1.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0523'
$ws.Range("E48").Value = '  +2.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.66'
$ws.Range("E49").Value = '  +0.09%  '
$ws.Range("E50").Value = '  +1.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0907'
$ws.Range("E51").Value = '  +1.54%  '
